$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1782.5
$ws.Range("J32").Value = 2512.2
$ws.Range("L32").Value = 2512.2
$ws.Range("N32").Value = -3164.2

$ws.Range("H98").Value = 2134.5557
$ws.Range("I98").Value = 1773.4286
$ws.Range("K98").Value = 1773.4286
$ws.Range("M98").Value = -275.4286

$ws.Range("H106").Value = 2320.9167
$ws.Range("I106").Value = 2320.9167
$ws.Range("K106").Value = 2320.9167
$ws.Range("M106").Value = -1689.9167

$ws.Range("H116").Value = 22837
$ws.Range("I116").Value = 100000
$ws.Range("K116").Value = 100000
$ws.Range("M116").Value = -96558

$ws.Range("H122").Value = 2134.5557
$ws.Range("I122").Value = 1773.4286
$ws.Range("K122").Value = 5320.2858
$ws.Range("M122").Value = -2870.2858

$ws.Range("H129").Value = 924.7538500000001
$ws.Range("J129").Value = 891.9508
$ws.Range("L129").Value = 2675.8524
$ws.Range("N129").Value = -12675.8524

$ws.Range("H131").Value = 2063.5293
$ws.Range("I131").Value = 697
$ws.Range("J131").Value = 4015.7144
$ws.Range("K131").Value = 2091
$ws.Range("L131").Value = 12047.1432
$ws.Range("M131").Value = 2949
$ws.Range("N131").Value = -22127.1432

$ws.Range("H137").Value = 1677
$ws.Range("I137").Value = 1425.25
$ws.Range("J137").Value = 1788.8889
$ws.Range("K137").Value = 4275.75
$ws.Range("L137").Value = 5366.6667
$ws.Range("M137").Value = -1725.75
$ws.Range("N137").Value = -10466.6667

$ws.Range("H141").Value = 3365.9
$ws.Range("I141").Value = 2379.8572
$ws.Range("K141").Value = 7139.571599999999
$ws.Range("M141").Value = -1959.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2486.4265
$ws.Range("I32").Value = 1642.5964
$ws.Range("J32").Value = 6859
$ws.Range("K32").Value = 1642.5964
$ws.Range("L32").Value = 6859
$ws.Range("M32").Value = -1355.5964
$ws.Range("N32").Value = -7433

$ws.Range("H74").Value = 5429.3335
$ws.Range("I74").Value = 5429.3335
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5429.3335
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4555.3335
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 5429.3335
$ws.Range("I77").Value = 5429.3335
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 27146.6675
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -22778.6675
$ws.Range("N77").ClearContents()

$ws.Range("H97").Value = 519.625
$ws.Range("I97").Value = 519.625
$ws.Range("K97").Value = 519.625
$ws.Range("M97").Value = -23.625

$ws.Range("H102").Value = 1250
$ws.Range("I102").Value = 1250
$ws.Range("K102").Value = 1250
$ws.Range("M102").Value = 372

$ws.Range("H109").Value = 67764.375
$ws.Range("J109").Value = 67764.375
$ws.Range("L109").Value = 67764.375
$ws.Range("N109").Value = -70538.375

$ws.Range("H122").Value = 1599.9
$ws.Range("I122").Value = 1576.9584
$ws.Range("K122").Value = 4730.8752
$ws.Range("M122").Value = -2280.8752

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1369.8
$ws.Range("I99").Value = 1212.25
$ws.Range("K99").Value = 1212.25
$ws.Range("M99").Value = 285.75

$ws.Range("H105").Value = 2371.875
$ws.Range("I105").Value = 2370.652
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 2370.652
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -623.652
$ws.Range("N105").Value = -5894

$ws.Range("H107").Value = 2240.8572
$ws.Range("I107").Value = 1962.125
$ws.Range("K107").Value = 1962.125
$ws.Range("M107").Value = -42.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2477.6316
$ws.Range("I31").Value = 1139.2693
$ws.Range("J31").Value = 5377.4165
$ws.Range("K31").Value = 1139.2693
$ws.Range("L31").Value = 5377.4165
$ws.Range("M31").Value = -844.2692999999999
$ws.Range("N31").Value = -5967.4165

$ws.Range("H34").Value = 2477.6316
$ws.Range("I34").Value = 1139.2693
$ws.Range("J34").Value = 5377.4165
$ws.Range("K34").Value = 1139.2693
$ws.Range("L34").Value = 5377.4165
$ws.Range("M34").Value = -937.2692999999999
$ws.Range("N34").Value = -5781.4165

$ws.Range("H58").Value = 1443.7142
$ws.Range("I58").Value = 1227
$ws.Range("J58").Value = 1606.25
$ws.Range("K58").Value = 1227
$ws.Range("L58").Value = 1606.25
$ws.Range("M58").Value = -1024
$ws.Range("N58").Value = -2012.25

$ws.Range("H62").Value = 4849.75
$ws.Range("I62").Value = 4499.5
$ws.Range("K62").Value = 4499.5
$ws.Range("M62").Value = -3875.5

$ws.Range("H65").Value = 4849.75
$ws.Range("I65").Value = 4499.5
$ws.Range("K65").Value = 22497.5
$ws.Range("M65").Value = -19377.5

$ws.Range("H99").Value = 2195
$ws.Range("I99").Value = 1660
$ws.Range("J99").Value = 2997.5
$ws.Range("K99").Value = 1660
$ws.Range("L99").Value = 2997.5
$ws.Range("M99").Value = -162
$ws.Range("N99").Value = -5993.5

$ws.Range("H126").Value = 2195
$ws.Range("I126").Value = 1660
$ws.Range("J126").Value = 2997.5
$ws.Range("K126").Value = 4980
$ws.Range("L126").Value = 8992.5
$ws.Range("M126").Value = -2510
$ws.Range("N126").Value = -13932.5

$ws.Range("H136").Value = 1443.7142
$ws.Range("I136").Value = 1227
$ws.Range("J136").Value = 1606.25
$ws.Range("K136").Value = 3681
$ws.Range("L136").Value = 4818.75
$ws.Range("M136").Value = -1131
$ws.Range("N136").Value = -9918.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 556
$ws.Range("I5").Value = 498.75
$ws.Range("J5").Value = 899.5
$ws.Range("K5").Value = 1496.25
$ws.Range("L5").Value = 2698.5
$ws.Range("M5").Value = -1384.25
$ws.Range("N5").Value = -2922.5

$ws.Range("H122").Value = 813.4
$ws.Range("I122").Value = 621.25
$ws.Range("J122").Value = 1033
$ws.Range("K122").Value = 5591.25
$ws.Range("L122").Value = 9297
$ws.Range("M122").Value = -3141.25
$ws.Range("N122").Value = -14197

$ws.Range("H135").Value = 556
$ws.Range("I135").Value = 498.75
$ws.Range("J135").Value = 899.5
$ws.Range("K135").Value = 4488.75
$ws.Range("L135").Value = 8095.5
$ws.Range("M135").Value = -1953.75
$ws.Range("N135").Value = -13165.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1321.4117
$ws.Range("J113").Value = 1441.7273
$ws.Range("L113").Value = 1441.7273
$ws.Range("N113").Value = -5781.7273

$ws.Range("H132").Value = 4040.077
$ws.Range("I132").Value = 3501
$ws.Range("J132").Value = 5253
$ws.Range("K132").Value = 10503
$ws.Range("L132").Value = 15759
$ws.Range("M132").Value = -7973
$ws.Range("N132").Value = -20819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 10805.4
$ws.Range("J43").Value = 10805.4
$ws.Range("L43").Value = 10805.4
$ws.Range("N43").Value = -11191.4

$ws.Range("H82").Value = 3173.7144
$ws.Range("I82").Value = 1790.4
$ws.Range("K82").Value = 1790.4
$ws.Range("M82").Value = -1429.4

$ws.Range("H85").Value = 3173.7144
$ws.Range("I85").Value = 1790.4
$ws.Range("K85").Value = 1790.4
$ws.Range("M85").Value = -542.4000000000001

$ws.Range("H136").Value = 4842.857
$ws.Range("I136").Value = 3685.7144
$ws.Range("K136").Value = 11057.1432
$ws.Range("M136").Value = -8507.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 160.83333
$ws.Range("I100").Value = 160.83333
$ws.Range("K100").Value = 321.66666
$ws.Range("M100").Value = 219.33334

$ws.Range("H107").Value = 1049.75
$ws.Range("I107").Value = 914
$ws.Range("K107").Value = 2742
$ws.Range("M107").Value = -822

$ws.Range("H113").Value = 762.625
$ws.Range("J113").Value = 920.2
$ws.Range("L113").Value = 2760.6
$ws.Range("N113").Value = -7100.6

$ws.Range("H126").Value = 3987.879
$ws.Range("J126").Value = 5441.5454
$ws.Range("L126").Value = 16324.6362
$ws.Range("N126").Value = -21264.6362

$ws.Range("H132").Value = 3620.1
$ws.Range("I132").Value = 1291.9166
$ws.Range("K132").Value = 3875.7498
$ws.Range("M132").Value = -1345.7498

$ws.Range("H136").Value = 4483.9473
$ws.Range("I136").Value = 4698.3335
$ws.Range("K136").Value = 14095.0005
$ws.Range("M136").Value = -11545.0005
